$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Footer timestamp update ---
$ws.Range("A1").Value = "Datos actualizados a 23 de Marzo de 2020 a las 23:46"

# --- Simple stat refreshes (no reordering) ---
# Estados Unidos (row 6)
$ws.Range("B6").Value = 43449
$ws.Range("C6").Value = 9883
$ws.Range("E6").Value = 42609
$ws.Range("G6").Value = 132
$ws.Range("H6").Value = 545

# Suiza (row 12)
$ws.Range("B12").Value = 8795
$ws.Range("C12").Value = 1321
$ws.Range("E12").Value = 8544
$ws.Range("G12").Value = 22
$ws.Range("H12").Value = 120

# Brasil (row 21)
$ws.Range("B21").Value = 1924
$ws.Range("C21").Value = 378
$ws.Range("E21").Value = 1888

# Australia (row 22)
$ws.Range("B22").Value = 1887
$ws.Range("C22").Value = 278
$ws.Range("E22").Value = 1762

# --- Nigeria re-ranked: moves up to directly after Afganistan (row 112) ---
# Rows 112-114 (Ruanda, Consejo Danes para los Refugiados, Mauricio) shift down one row,
# and Nigeria (previously row 115) takes row 112 with updated figures.
$ws.Range("A115").Value = "Mauricio"
$ws.Range("B115").Value = 36
$ws.Range("C115").Value = 8
$ws.Range("D115").Value = 0
$ws.Range("E115").Value = 34
$ws.Range("F115").Value = 1
$ws.Range("G115").Value = 0
$ws.Range("H115").Value = 2

$ws.Range("A114").Value = "Consejo Danes para los Refugiados"
$ws.Range("B114").Value = 36
$ws.Range("C114").Value = 6
$ws.Range("D114").Value = 0
$ws.Range("E114").Value = 35
$ws.Range("F114").Value = 0
$ws.Range("G114").Value = 0
$ws.Range("H114").Value = 1

$ws.Range("A113").Value = "Ruanda"
$ws.Range("B113").Value = 36
$ws.Range("C113").Value = 17
$ws.Range("D113").Value = 0
$ws.Range("E113").Value = 36
$ws.Range("F113").Value = 0
$ws.Range("G113").Value = 0
$ws.Range("H113").Value = 0

$ws.Range("A112").Value = "Nigeria"
$ws.Range("B112").Value = 40
$ws.Range("C112").Value = 10
$ws.Range("D112").Value = 2
$ws.Range("E112").Value = 37
$ws.Range("F112").Value = 0
$ws.Range("G112").Value = 1
$ws.Range("H112").Value = 1
